$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 31: new reference entry (Hasan/Kabir book)
$ws.Range("A31").Value = "Hasan A.R., Kabir C.S., ""Fluid Flow and Heat Transfer in Wellbores"", SPE, Richardson, Texas, 2002"
$ws.Range("B31").Value = "часто цитируемая книга, пока не нашел в открытом доступе"

# Row 32: new reference entry (Кислицын/Шабаров book)
$ws.Range("A32").Value = "Кислицын А. А. Тепломасообмен / А. А. Кислицын, А. Б. Шабаров. Тюмень: `nизд-во ТюмГУ, 2008. "

# Row 33: new reference entry (Вакулин/Шабаров book)
$ws.Range("A33").Value = "Вакулин А. А. Диагностика теплофизических параметров в нефтегазовых `nтехнологиях / А. А. Вакулин, А. Б. Шабаров. Новосибирск: Наука. `nСиб. Издательская фирма РАН, 1998."

# D27: new note about hydrate articles
$ws.Range("D27").Value = "статьи по гидратам стр.6, 22, 50, 84, 92, 102"

# D16: quote-prefixed "++'" text (leading apostrophe marks it as text-with-prefix)
$ws.Range("D16").Value = "'++'"

# Row heights
$ws.Rows.Item(31).RowHeight = 52.9
$ws.Rows.Item(32).RowHeight = 42.75
$ws.Rows.Item(33).RowHeight = 85.5

# Styling: D16 uses quote-prefix style (font 0, border 0, quotePrefix)
$ws.Range("D16").HorizontalAlignment = -4108  # xlCenter (default-ish, no border needed since quotePrefix style has borderId=0)

# A31: new font (Times New Roman 14) with justify/center alignment
$fontA31 = $ws.Range("A31").Font
$fontA31.Name = "Times New Roman"
$fontA31.Size = 14
$ws.Range("A31").HorizontalAlignment = -4130  # xlJustify
$ws.Range("A31").VerticalAlignment = -4108  # xlCenter

# B31: bordered box style, centered, wraptext
$ws.Range("B31").Borders.LineStyle = 1
$ws.Range("B31").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B31").VerticalAlignment = -4108  # xlCenter
$ws.Range("B31").WrapText = $true

# A32/A33 reuse the same border+center+wrap format already used by column B link cells (e.g. B4)
$ws.Range("B4").Copy()
$ws.Range("A32").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A33").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update dimension / view state is handled automatically by Excel on save.
$ws.Application.ActiveWindow.ScrollRow = 26
$ws.Range("D31").Select()
